$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 12.73493403024474
$ws.Range("C2").Value = 9.341912471161857
$ws.Range("E2").Value = 11.59650481190773
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 3.602027665198992
$ws.Range("M2").Value = 14.48695727445525
$ws.Range("O2").Value = 17.82169455400742
$ws.Range("B3").Value = 12.07281554817309
$ws.Range("C3").Value = 8.919560778995537
$ws.Range("E3").Value = 11.52645985979907
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 3.604208836831901
$ws.Range("M3").Value = 14.16670987279553
$ws.Range("O3").Value = 17.96004501333924
$ws.Range("B4").Value = 11.64735031185373
$ws.Range("C4").Value = 8.648664134182619
$ws.Range("E4").Value = 11.48852497756657
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 3.605617137558023
$ws.Range("M4").Value = 13.96911603980379
$ws.Range("O4").Value = 18.05209657827689
$ws.Range("B5").Value = 11.46938119765977
$ws.Range("C5").Value = 8.535466857338998
$ws.Range("E5").Value = 11.4743536280604
$ws.Range("F5").Value = 15.008197319934
$ws.Range("G5").Value = 3.60620845323629
$ws.Range("M5").Value = 13.88846814038596
$ws.Range("O5").Value = 18.09138340818261
$ws.Range("B6").Value = 11.43955799623797
$ws.Range("C6").Value = 8.516504531842267
$ws.Range("E6").Value = 11.47207853141475
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("G6").Value = 3.606307694683303
$ws.Range("M6").Value = 13.87507238016601
$ws.Range("O6").Value = 18.09801384050647
$ws.Range("B7").Value = 11.64496849020826
$ws.Range("C7").Value = 8.647148721440065
$ws.Range("E7").Value = 11.48832863164641
$ws.Range("F7").Value = 15.26647399323133
$ws.Range("G7").Value = 3.605625041631954
$ws.Range("M7").Value = 13.96802875330765
$ws.Range("O7").Value = 18.0526192418134
$ws.Range("B8").Value = 12.51065514380649
$ws.Range("C8").Value = 9.198745151714244
$ws.Range("E8").Value = 11.57130834745978
$ws.Range("F8").Value = 16.5399640634477
$ws.Range("G8").Value = 3.602765436251643
$ws.Range("M8").Value = 14.37680305569937
$ws.Range("O8").Value = 17.86791658999276
$ws.Range("B9").Value = 14.05212923173297
$ws.Range("C9").Value = 10.18478481584275
$ws.Range("E9").Value = 11.77362792893048
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 3.597703010243848
$ws.Range("M9").Value = 15.16594648511836
$ws.Range("O9").Value = 17.56258997138859
$ws.Range("B10").Value = 15.08312123276345
$ws.Range("C10").Value = 10.84674326838678
$ws.Range("E10").Value = 11.9453462303883
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 3.59431232533761
$ws.Range("M10").Value = 15.73199131910986
$ws.Range("O10").Value = 17.37364497244184
$ws.Range("B11").Value = 15.52915474062648
$ws.Range("C11").Value = 11.13366300713865
$ws.Range("E11").Value = 12.02819788305429
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 3.592840391007122
$ws.Range("M11").Value = 15.98532714841199
$ws.Range("O11").Value = 17.29552595441872
$ws.Range("B12").Value = 15.69469410363464
$ws.Range("C12").Value = 11.24022642424641
$ws.Range("E12").Value = 12.06022628350726
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 3.59229308610443
$ws.Range("M12").Value = 16.08056961442065
$ws.Range("O12").Value = 17.26708378992211
$ws.Range("B13").Value = 15.6591926837949
$ws.Range("C13").Value = 11.21736951256154
$ws.Range("E13").Value = 12.05329974038984
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 3.592410510343055
$ws.Range("M13").Value = 16.06008954389522
$ws.Range("O13").Value = 17.27315841012333
$ws.Range("B14").Value = 15.54284148396159
$ws.Range("C14").Value = 11.14247208994231
$ws.Range("E14").Value = 12.03081992582532
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 3.592795162098619
$ws.Range("M14").Value = 15.99317714253003
$ws.Range("O14").Value = 17.29316307946312
$ws.Range("B15").Value = 15.47113335272292
$ws.Range("C15").Value = 11.09632232219911
$ws.Range("E15").Value = 12.01713477136864
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 3.593032084135693
$ws.Range("M15").Value = 15.95209884420472
$ws.Range("O15").Value = 17.30556535529002
$ws.Range("B16").Value = 15.05350418987891
$ws.Range("C16").Value = 10.82770255970685
$ws.Range("E16").Value = 11.94002473641378
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 3.59440993383925
$ws.Range("M16").Value = 15.7153434910567
$ws.Range("O16").Value = 17.37890899947833
$ws.Range("B17").Value = 14.79137166711551
$ws.Range("C17").Value = 10.65924023479195
$ws.Range("E17").Value = 11.89391480806212
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 3.595273219642264
$ws.Range("M17").Value = 15.5689685914175
$ws.Range("O17").Value = 17.42591868611324
$ws.Range("B18").Value = 14.63844067136333
$ws.Range("C18").Value = 10.56100985282031
$ws.Range("E18").Value = 11.86784111255613
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 3.595776398287347
$ws.Range("M18").Value = 15.48439108622179
$ws.Range("O18").Value = 17.45369388157662
$ws.Range("B19").Value = 14.58629206651637
$ws.Range("C19").Value = 10.52752285622135
$ws.Range("E19").Value = 11.85909065399687
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 3.595947907916867
$ws.Range("M19").Value = 15.45569111614201
$ws.Range("O19").Value = 17.46322416742918
$ws.Range("B20").Value = 14.81950008198776
$ws.Range("C20").Value = 10.67731187795455
$ws.Range("E20").Value = 11.89877716060595
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 3.595180634623812
$ws.Range("M20").Value = 15.58459114089139
$ws.Range("O20").Value = 17.42083810620744
$ws.Range("B21").Value = 15.57710836812272
$ws.Range("C21").Value = 11.16452822484586
$ws.Range("E21").Value = 12.03740525757947
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 3.5926819072866
$ws.Range("M21").Value = 16.01285037315263
$ws.Range("O21").Value = 17.28725618124058
$ws.Range("B22").Value = 16.05262231484438
$ws.Range("C22").Value = 11.47077346847118
$ws.Range("E22").Value = 12.13180674388927
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 3.591107601009016
$ws.Range("M22").Value = 16.28868071133488
$ws.Range("O22").Value = 17.20660269402184
$ws.Range("B23").Value = 15.80064437621179
$ws.Range("C23").Value = 11.3084514461858
$ws.Range("E23").Value = 12.08108448553497
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 3.591942479393298
$ws.Range("M23").Value = 16.14186526260201
$ws.Range("O23").Value = 17.24903602787394
$ws.Range("B24").Value = 14.80679015912716
$ws.Range("C24").Value = 10.66914598049675
$ws.Range("E24").Value = 11.89657753242009
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 3.59522247092439
$ws.Range("M24").Value = 15.57752950182717
$ws.Range("O24").Value = 17.42313270516409
$ws.Range("B25").Value = 13.65257189137662
$ws.Range("C25").Value = 9.928752274866834
$ws.Range("E25").Value = 11.71475604060269
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 3.59901454351951
$ws.Range("M25").Value = 14.95445135508748
$ws.Range("O25").Value = 17.63902782226981
